$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Optimizer | 1000 run - Descend")

$ws.Range("C3").Value = 281733
$ws.Range("D3").Value = 320296
$ws.Range("I3").Value = 62254
$ws.Range("C4").Value = 281733
$ws.Range("D4").Value = 924358
$ws.Range("I4").Value = 179461
$ws.Range("C5").Value = 281733
$ws.Range("D5").Value = 1589268
$ws.Range("I5").Value = 338096
$ws.Range("C6").Value = 281733
$ws.Range("D6").Value = 2247578
$ws.Range("I6").Value = 526631
$ws.Range("C7").Value = 281733
$ws.Range("D7").Value = 2834414
$ws.Range("I7").Value = 721877
$ws.Range("C8").Value = 281733
$ws.Range("D8").Value = 2899288
$ws.Range("I8").Value = 745066
$ws.Range("C9").Value = 281733
$ws.Range("D9").Value = 2964096
$ws.Range("I9").Value = 768554
$ws.Range("C10").Value = 281733
$ws.Range("D10").Value = 3028838
$ws.Range("I10").Value = 792341
$ws.Range("C11").Value = 281733
$ws.Range("D11").Value = 3093514
$ws.Range("I11").Value = 816427
$ws.Range("C12").Value = 281733
$ws.Range("D12").Value = 3158124
$ws.Range("I12").Value = 840812
$ws.Range("C13").Value = 281733
$ws.Range("D13").Value = 3222668
$ws.Range("I13").Value = 865496
$ws.Range("C14").Value = 281733
$ws.Range("D14").Value = 3287146
$ws.Range("I14").Value = 890479
$ws.Range("C15").Value = 281733
$ws.Range("D15").Value = 3351558
$ws.Range("I15").Value = 915761
$ws.Range("C16").Value = 281733
$ws.Range("D16").Value = 3415904
$ws.Range("I16").Value = 941342
$ws.Range("C17").Value = 281733
$ws.Range("D17").Value = 3480184
$ws.Range("I17").Value = 967222
$ws.Range("C18").Value = 281733
$ws.Range("D18").Value = 3544398
$ws.Range("I18").Value = 993401
$ws.Range("C19").Value = 281733
$ws.Range("D19").Value = 4182908
$ws.Range("I19").Value = 1271636
$ws.Range("C20").Value = 281733
$ws.Range("D20").Value = 4814818
$ws.Range("I20").Value = 1579771
$ws.Range("C21").Value = 281733
$ws.Range("D21").Value = 5440128
$ws.Range("I21").Value = 1917806
$ws.Range("C22").Value = 281733
$ws.Range("D22").Value = 6058838
$ws.Range("I22").Value = 2285741
$ws.Range("C23").Value = 281733
$ws.Range("D23").Value = 6653624
$ws.Range("I23").Value = 2678776

$ws.Activate()

